# mantenimiento inspeccion y parte de inventario
#
# Adds a "LICENCIA DE TRANSITO" row (with plate/license numbers) to the
# GTK801 and JXV805 equipment sheets, and renames GTK801's "NUMERO DE
# REGISTRO" label to "PLACA".

$wb = $excel.ActiveWorkbook

# xlEdge* constants used for per-edge border assignment.
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlCenter = -4108

function Add-LicenciaRow {
    param($ws, $RowIndex, $PlacaValue)

    # Push the existing rows (NUMERO DE MOTOR / NUMERO DE SERIE / PROPIETARIO, ...)
    # down by one to make room for the new row, mirroring an in-UI "Insert Row".
    $ws.Rows.Item($RowIndex).Insert()

    $leftRange = $ws.Range("A" + $RowIndex + ":B" + $RowIndex)
    $leftRange.Merge()
    $leftRange.Value = "LICENCIA DE TRANSITO"
    $leftRange.HorizontalAlignment = $xlCenter

    $rightRange = $ws.Range("C" + $RowIndex + ":D" + $RowIndex)
    $rightRange.Merge()
    $rightRange.Value = $PlacaValue
    $rightRange.HorizontalAlignment = $xlCenter

    # Give the merged pair a continuous-looking box: the left cell of each
    # merge gets left+top+bottom, the right cell gets right+top+bottom, so
    # there's no divider line down the middle of the merged range.
    $a = $ws.Cells.Item($RowIndex, 1)
    $a.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $a.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $a.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

    $b = $ws.Cells.Item($RowIndex, 2)
    $b.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $b.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $b.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

    $c = $ws.Cells.Item($RowIndex, 3)
    $c.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

    $d = $ws.Cells.Item($RowIndex, 4)
    $d.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $d.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $d.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
}

# ----- GTK801 -----
$wsGtk = $wb.Worksheets.Item("GTK801")

# "NUMERO DE REGISTRO" -> "PLACA"
$wsGtk.Range("A14").Value = "PLACA"

Add-LicenciaRow $wsGtk 16 10020357127

$wsGtk.Activate()
$wsGtk.Range("A16:D16").Select()

# ----- JXV805 -----
$wsJxv = $wb.Worksheets.Item("JXV805")

Add-LicenciaRow $wsJxv 16 10027092604

$wsJxv.Activate()
$wsJxv.Range("C16:D16").Select()
